$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.1124019449155688, 0.6821200353019861, 1.209983865900432, 1.099992666293931, 1.105344008082418, 50),
    @(0.1960116252932615, 0.6752669640606744, 1.511496450551623, 1.22942931905483, 1.226547159135974, 48),
    @(0.001728000735583749, 0.6309824882362022, 1.099648267891479, 1.048641153060225, 1.059736755062583, 48),
    @(0.0937284794034382, 0.6894373285789633, 1.642876430466323, 1.281747412896286, 1.292135878597467, 47),
    @(0.04836703510120018, 0.6561079152665072, 1.680400965917633, 1.296302806414316, 1.309714416778838, 46),
    @(0.0481232454536701, 0.6637583974150937, 1.194123607785736, 1.092759629463742, 1.108116937751023, 34),
    @(0.09097433578834678, 0.7080043655050214, 1.193680625051876, 1.092556920737714, 1.105643785089591, 33),
    @(0.03775263814022831, 0.6203147181500017, 0.600319873274656, 0.7748031190403508, 0.7992627346669655, 16),
    @(0.02639662487600095, 0.6628602414482503, 0.693693318076211, 0.8328825355812253, 0.8774942458002957, 10),
    @(0.2236099196487949, 0.5232375926460167, 0.3429526541974182, 0.5856215964233373, 0.6051355819484564, 5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
